{"js": "// Split the run containing \"{{ juego }}\" into two runs:\n//   \"{{\"        (unchanged formatting)\n//   \" juego }}\" (unchanged formatting)\n// without altering any other run in the paragraph.\n//\n// A plain text-insert (Range.insertText) on this engine rebuilds the whole\n// paragraph and coalesces every adjacent run that shares identical\n// formatting, which would wreck the surrounding runs. Inserting (and then\n// removing) a bookmark at the split point forces a clean, local run split\n// instead, leaving every other run untouched.\n\nconst target = \"{{ juego }}\";\nconst results = context.document.body.search(target, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"' + target + '\" in the document.');\n}\n\nconst fullRange = results.items[0];\n\n// Locate the boundary right after \"{{\" inside the matched range.\nconst braceResults = fullRange.search(\"{{\", { matchCase: true });\nbraceResults.load(\"items\");\nawait context.sync();\n\nconst splitPoint = braceResults.items[0].getRange(\"After\");\n\nconst bookmarkName = \"__splitJuegoTmp\";\nsplitPoint.insertBookmark(bookmarkName);\nawait context.sync();\n\ncontext.document.deleteBookmark(bookmarkName);\nawait context.sync();\n", "ps1": "# Split the run containing \"{{ juego }}\" into two runs:\n#   \"{{\"        (unchanged formatting)\n#   \" juego }}\" (unchanged formatting)\n# without altering any other run in the paragraph.\n#\n# A plain text replace (Range.Text = ...) on this engine rebuilds the whole\n# paragraph and coalesces every adjacent run that shares identical\n# formatting, which would wreck the surrounding runs (\"(\", \"s\", \")\", etc.).\n# Inserting (and then removing) a temporary bookmark at the split point\n# forces a clean, local run split instead, leaving every other run intact.\n\n$d = $word.ActiveDocument\n\n$target = $d.Content\n$found = $target.Find.Execute(\"{{ juego }}\")\n\nif ($found) {\n    # $target now spans exactly \"{{ juego }}\"; duplicate it and search\n    # inside for \"{{\" to find the precise split boundary.\n    $braceRng = $target.Duplicate\n    $braceFound = $braceRng.Find.Execute(\"{{\")\n\n    if ($braceFound) {\n        $splitPoint = $d.Range($braceRng.End, $braceRng.End)\n\n        $bookmarkName = \"__splitJuegoTmp\"\n        $d.Bookmarks.Add($bookmarkName, $splitPoint)\n        $d.Bookmarks($bookmarkName).Delete()\n    }\n}\n"}
